$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45; this pushes the existing rows 45-138 down to 46-139
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 45152
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112038
$ws.Range("G45").Value = "Cebollín baby"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = 1467
$ws.Range("N45").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 734
$ws.Range("Q45").Value = 2
$ws.Range("R45").Value = "Hortaliza"
